$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 5 values to reduced precision (custom accuracy to 2 decimals)
$ws.Range("B5").Value = 18.74
$ws.Range("C5").Value = 13.73
$ws.Range("D5").Value = 1.18
$ws.Range("E5").Value = 40.69
$ws.Range("F5").Value = 33.28
$ws.Range("G5").Value = 14.74
$ws.Range("H5").Value = 53.4
$ws.Range("I5").Value = 22.69
$ws.Range("J5").Value = 10.02
$ws.Range("K5").Value = 14.84
$ws.Range("L5").Value = 16.34
$ws.Range("M5").Value = 17.19
$ws.Range("N5").Value = 4.71
$ws.Range("O5").Value = 14.66
$ws.Range("P5").Value = 20.81
$ws.Range("Q5").Value = 12.41
$ws.Range("R5").Value = 0.83
$ws.Range("S5").Value = 0.78
$ws.Range("T5").Value = 215.84
$ws.Range("U5").Value = 40.89
$ws.Range("V5").Value = 13.53
$ws.Range("W5").Value = 27.43
$ws.Range("X5").Value = 14.4
$ws.Range("Y5").Value = 2.23
$ws.Range("Z5").Value = 26.38
$ws.Range("AA5").Value = 11.95
$ws.Range("AB5").Value = 10.64
$ws.Range("AC5").Value = 12.5
$ws.Range("AD5").Value = 17.09
$ws.Range("AE5").Value = 0.56
$ws.Range("AG5").Value = 7.6

# Remove the last data row (row 6) entirely - trimming the dataset
$ws.Rows(6).Delete()

# Narrow column J (10th column) from width 8 to width 7
$ws.Columns(10).ColumnWidth = 6.17
